$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.653.56'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').Value = '3.097.37'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.98'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '615.52'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.08'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -5.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.391'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +7.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '3.094.60'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('E11').Value = '  -0.82%  '
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000251'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.92%  '
$ws.Range('D14').Value = '91.963.43'
$ws.Range('E14').Value = '  +2.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.52'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.04'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('D17').Value = '3.674.51'
$ws.Range('E17').Value = '  +0.30%  '
$ws.Range('D18').Value = '3.108.00'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('E19').Value = '  -1.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.71'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.82'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.40%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.30'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.68%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '442.77'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000200'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.54%  '
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '86.54'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.63'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.32%  '
$ws.Range('D28').Value = '3.270.56'
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.137'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +15.10%  '
$ws.Range('E31').Value = '  -7.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.167'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -5.20%  '
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('E34').Value = '  +3.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.169'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.91'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('E38').Value = '  -5.46%  '
$ws.Range('E39').Value = '  +0.96%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.30'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.30%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '477.97'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.68%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.431'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.46%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.37'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -5.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.18'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '158.93'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.37%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.90'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.31%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.696'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.80%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0341'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +9.70%  '
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.36'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.92'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.67%  '
